$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.478339791297913
$ws.Range("B1").Value = 4.166070461273193
$ws.Range("C1").Value = 3.463298559188843
$ws.Range("D1").Value = 1.889957904815674
$ws.Range("E1").Value = 0.6581440567970276
